# EnrollTaskData.xlsx: hide inactive tasks (deselect the old active "學生"
# sheet, resetting its cursor) and add a new "禁止選課名單" (forbidden
# course-selection list) sheet of students who are barred from enrolling.

$wb = $excel.ActiveWorkbook

# The "學生" (Students) sheet was the active tab with the cursor parked at
# D10; it is no longer the active sheet, so its selection resets to A1.
$wsStudents = $wb.Worksheets.Item("學生")
$wsStudents.Range("A1").Select()

# Add the new sheet after all existing sheets (it becomes the active tab).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForbidden = $wb.Worksheets.Add($null, $lastSheet)
$wsForbidden.Name = "禁止選課名單"

# Header row: same student-identifying columns as "學生", plus a new
# "禁止原因" (reason for prohibition) column.
$wsForbidden.Range("A1").Value = "學號"
$wsForbidden.Range("B1").Value = "班級"
$wsForbidden.Range("C1").Value = "座號"
$wsForbidden.Range("D1").Value = "姓名"
$wsForbidden.Range("E1").Value = "禁止原因"

# Leave the cursor on the new sheet at E2, ready for data entry.
$wsForbidden.Range("E2").Select()
